# Add a new "options" dialog block (options/music/sound/speech/on/off/close
# key-value pairs) to the Language sheet, inserted right after the existing
# header rows (welcome/title/test/test2) and before the tag_/material_ rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows starting at row 6 - this shifts all the existing
# tag_/material_ rows down by 7 (old row 6 -> new row 13, etc.) while Excel
# keeps every shared-string reference intact automatically.
$ws.Rows("6:12").Insert()

# Fill in the newly-inserted rows with the new localization keys.
$ws.Range("A6").Value = "options"
$ws.Range("B6").Value = "OPTIONS"

$ws.Range("A7").Value = "music"
$ws.Range("B7").Value = "MUSIC"

$ws.Range("A8").Value = "sound"
$ws.Range("B8").Value = "SOUND"

$ws.Range("A9").Value = "speech"
$ws.Range("B9").Value = "SPEECH"

$ws.Range("A10").Value = "on"
$ws.Range("B10").Value = "ON"

$ws.Range("A11").Value = "off"
$ws.Range("B11").Value = "OFF"

$ws.Range("A12").Value = "close"
$ws.Range("B12").Value = "CLOSE"

# Match the author's resulting view state: scrolled back to the top, with
# the new "close" key row selected.
$ws.Range("A12").Select()
